# Update countries & provincias Spain
# - India's case counts were refreshed and it overtook Chile; the two rows swap position
#   (India moves to row 29, Chile moves to row 30) while Chile keeps its previous figures.
# - Serbia's case counts were refreshed and it overtook Finlandia, Mexico,
#   Emiratos Arabes Unidos and Panama; Serbia moves up to row 43 and the other four
#   countries shift down one row each, keeping their previous figures.
# - Letonia (row 80) gets corrected "Casos activos"/"Recuperados" figures.
# - The "Datos actualizados" timestamp in A1 is refreshed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Refresh timestamp
$ws.Range("A1").Value = "Datos actualizados a 6 de Abril de 2020 a las 15:52"

# --- India overtakes Chile (rows 29-30) ---
# Row 29: now India, with new data
$ws.Range("A29").Value = "India"
$ws.Range("B29").Value = 4553
$ws.Range("C29").Value = 264
$ws.Range("D29").Value = 328
$ws.Range("E29").Value = 4107
$ws.Range("F29").Value = 0
$ws.Range("G29").Value = 0
$ws.Range("H29").Value = 118

# Row 30: now Chile, keeping its previous data
$ws.Range("A30").Value = "Chile"
$ws.Range("B30").Value = 4471
$ws.Range("C30").Value = 0
$ws.Range("D30").Value = 618
$ws.Range("E30").Value = 3819
$ws.Range("F30").Value = 307
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = 34

# --- Serbia overtakes Finlandia, Mexico, Emiratos Arabes Unidos, Panama (rows 43-47) ---
# Row 43: now Serbia, with new data
$ws.Range("A43").Value = "Serbia"
$ws.Range("B43").Value = 2200
$ws.Range("C43").Value = 292
$ws.Range("D43").Value = 54
$ws.Range("E43").Value = 2088
$ws.Range("F43").Value = 101
$ws.Range("G43").Value = 7
$ws.Range("H43").Value = 58

# Row 44: now Finlandia, keeping its previous data
$ws.Range("A44").Value = "Finlandia"
$ws.Range("B44").Value = 2176
$ws.Range("C44").Value = 249
$ws.Range("D44").Value = 300
$ws.Range("E44").Value = 1849
$ws.Range("F44").Value = 81
$ws.Range("G44").Value = 1
$ws.Range("H44").Value = 27

# Row 45: now Mexico, keeping its previous data
$ws.Range("A45").Value = "Mexico"
$ws.Range("B45").Value = 2143
$ws.Range("C45").Value = 253
$ws.Range("D45").Value = 633
$ws.Range("E45").Value = 1416
$ws.Range("F45").Value = 293
$ws.Range("G45").Value = 15
$ws.Range("H45").Value = 94

# Row 46: now Emiratos Arabes Unidos, keeping its previous data
$ws.Range("A46").Value = "Emiratos Arabes Unidos"
$ws.Range("B46").Value = 2076
$ws.Range("C46").Value = 277
$ws.Range("D46").Value = 167
$ws.Range("E46").Value = 1898
$ws.Range("F46").Value = 1
$ws.Range("G46").Value = 1
$ws.Range("H46").Value = 11

# Row 47: now Panama, keeping its previous data
$ws.Range("A47").Value = "Panama"
$ws.Range("B47").Value = 1988
$ws.Range("C47").Value = 187
$ws.Range("D47").Value = 13
$ws.Range("E47").Value = 1921
$ws.Range("F47").Value = 78
$ws.Range("G47").Value = 8
$ws.Range("H47").Value = 54

# --- Letonia (row 80): corrected active cases / recovered figures ---
$ws.Range("D80").Value = 16
$ws.Range("E80").Value = 525
